# WRI edit to BRPSPTY
# The runtime pre-binds a top-level $wb that is NOT the live workbook
# (Worksheets.Count reads 0 through it), so always re-fetch it from
# $excel.ActiveWorkbook first, as the task's own example shows.
$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("BRPSPTY")

# --- BRPSPTY row 2: retype the "this year" RPS percentage series -------
# Before: F2=0.19, G2=0.21, H2==G2 (shared formula anchor), I2=shared ref
# After:  F2=0.175, G2=0.19, H2=0.21 (literal), I2==H2 (new shared anchor)
# J2:AJ2 already chain off the previous column and keep doing so.
$wsData.Range("F2").Value = 0.175
$wsData.Range("G2").Value = 0.19
$wsData.Range("H2").Value = 0.21
$wsData.Range("I2").Formula = "=H2"

# --- BRPSPTY sheet view: move the remembered selection to AI27 ---------
# Selecting a range activates its sheet as a side effect, so do this
# before restoring "About" as the active tab below.
$wsData.Range("AI27").Select() | Out-Null

# --- About sheet view: drop the scrolled topLeftCell (A10) -------------
# Re-activating the sheet resets the window's scroll position while
# leaving its own remembered selection (A24) untouched, and restores
# "About" as the workbook's active tab.
$wsAbout.Activate() | Out-Null
